$wb = $excel.ActiveWorkbook

# --- Add the new "financial ratios" worksheet after the last existing sheet ---
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "financial ratios"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 18.3
$ws.Columns.Item(4).ColumnWidth = 18.0
$ws.Columns.Item(5).ColumnWidth = 18.5
$ws.Columns.Item(6).ColumnWidth = 18.0

# --- Row 5: section title ---
$ws.Range("B5").Value = "Coefficients of determination"
$ws.Range("B5").Font.Bold = $true

# --- Row 6: headers ---
$ws.Range("B6").Value = "Financial ratio"
$ws.Range("C6").Value = "Decision tree (train)"
$ws.Range("D6").Value = "Decision tree (test)"
$ws.Range("E6").Value = "Gradient boost (train)"
$ws.Range("F6").Value = "Gradient boost (test)"

# --- Data rows 7-14 ---
$labels = @("Enterprise value", "Free cash flow", "EBITDA", "Revenue", "Return on equity", "Gross profit margin", "Quick ratio", "Debt to equity ratio")
$c = @(1, 1, 1, 1, 1, 0.65048829487508, 0.611219829794956, 0.611219829794956)
$d = @(-2.58198612156547, -5.05831024318142, -8.47714095853272, -6.64862911206342, -6.64862911206342, -8.31246456204073, -6.64862911206342, -6.64862911206342)
$e = @(0.999999995658277, 0.999999994794889, 0.999999896606301, 0.999999993985728, 0.999999996845121, 0.422532456113729, 0.292273915800679, 0.293000182550571)
$f = @(-5.16351442656913, -7.7442853512372, -9.07728458445717, -5.56904110699226, -4.43375034207515, -4.96361618177758, -5.16208651074748, -9.32353157453062)

for ($i = 0; $i -lt 8; $i++) {
    $row = 7 + $i
    $ws.Range("B$row").Value = $labels[$i]
    $ws.Range("C$row").Value = $c[$i]
    $ws.Range("D$row").Value = $d[$i]
    $ws.Range("E$row").Value = $e[$i]
    $ws.Range("F$row").Value = $f[$i]
}

$ws.Range("C7:F14").NumberFormat = "0.000"

# --- View state: record the last selection on the new sheet ---
$ws.Range("E14").Select() | Out-Null

# --- Restore view state on the "sector" sheet (topLeftCell moves to A1) ---
$wsSector = $wb.Worksheets.Item("sector")
$wsSector.Activate() | Out-Null
$wsSector.Range("F69").Select() | Out-Null

# --- "macro" stays the active/selected tab, with a new active cell B6 ---
$wsMacro = $wb.Worksheets.Item("macro")
$wsMacro.Activate() | Out-Null
$wsMacro.Range("B6").Select() | Out-Null
